# Update the localized "age"/"education" section headers (rows 19 and 29)
# on the single worksheet so each now reads "By <category>" / "По <category>"
# instead of the bare category name, matching the re-worded headers used
# elsewhere in the workbook (e.g. "По территории" / "By territory").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A (Kyrgyz): "age (in years)" and "education" section headers
$ws.Range("A19").Value = "Жаш курагы боюнча (жылдарда)"
$ws.Range("A29").Value = "Билими боюнча"

# Column B (Russian): "age (in years)" and "education" section headers
$ws.Range("B19").Value = "По возрасту (в годах)"
$ws.Range("B29").Value = "По образованию"

# Column C (English): "age (in years)" and "education" section headers
$ws.Range("C19").Value = "By age (in years) "
$ws.Range("C29").Value = "By education"
